# Add a new "Sheet2" right after "Sheet1", make it the active sheet,
# give it a bold header row ("ANO" / "PHOTOS"), and fill in the
# year -> photo-count data.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($null, $sheet1)
$ws.Name = "Sheet2"

# Header row (bold)
$ws.Range("A1").Value = "ANO"
$ws.Range("B1").Value = "PHOTOS"
$ws.Range("A1:B1").Font.Bold = $true

# Data rows
$data = @(
    @(2011, 16),
    @(2014, 8318),
    @(2015, 27009),
    @(2016, 33789),
    @(2017, 3938),
    @(2018, 96022),
    @(2019, 512519),
    @(2020, 8539),
    @(2021, 1741455),
    @(2022, 20),
    @(2023, 1413),
    @(2024, 144)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row = $row + 1
}

$ws.Activate()
